# Re-shuffle the data rows (2-8) of the "Corazón de apio" sheet so that the
# record content ends up matching the new arrangement described by the diff.
# Row 6 stays the same; the others get new values for columns D, I, J, K, L,
# M, N, P (Q only changes along with row 4/8 swap).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 45225
$ws.Range("J2").Value = 60
$ws.Range("L2").Value = 2000
$ws.Range("M2").Value = 1750
$ws.Range("P2").Value = 292

# Row 3
$ws.Range("D3").Value = 44267
$ws.Range("J3").Value = 120
$ws.Range("L3").Value = 1800
$ws.Range("M3").Value = 1650
$ws.Range("P3").Value = 275

# Row 4
$ws.Range("D4").Value = 44623
$ws.Range("J4").Value = 300
$ws.Range("K4").Value = 1800
$ws.Range("L4").Value = 2000
$ws.Range("M4").Value = 1900
$ws.Range("N4").Value = "$/paquete"
$ws.Range("P4").Value = 1900
$ws.Range("Q4").Value = 1

# Row 5
$ws.Range("D5").Value = 44377
$ws.Range("I5").Value = "Segunda"
$ws.Range("J5").Value = 550
$ws.Range("K5").Value = 2000
$ws.Range("L5").Value = 2800
$ws.Range("M5").Value = 2364
$ws.Range("P5").Value = 394

# Row 7
$ws.Range("D7").Value = 45204
$ws.Range("I7").Value = "Primera"
$ws.Range("J7").Value = 100
$ws.Range("K7").Value = 800
$ws.Range("L7").Value = 1000
$ws.Range("M7").Value = 900
$ws.Range("P7").Value = 150

# Row 8
$ws.Range("D8").Value = 45218
$ws.Range("J8").Value = 180
$ws.Range("K8").Value = 1400
$ws.Range("L8").Value = 1500
$ws.Range("M8").Value = 1444
$ws.Range("N8").Value = "$/docena de matas"
$ws.Range("P8").Value = 241
$ws.Range("Q8").Value = 6
